$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Header for new column L
$ws.Range("L1").Value = "Clima"

# Values for column L, rows 2-27 (climate classification per appezzamento)
$climaValues = @(
    "Mediterraneo_di_transizione", # row 2
    "Temperato",                   # row 3
    "Mediterraneo_di_transizione", # row 4
    "Temperato",                   # row 5
    "Mediterraneo_di_transizione", # row 6
    "Mediterraneo_di_transizione", # row 7
    "Mediterraneo_di_transizione", # row 8
    "Mediterraneo_di_transizione", # row 9
    "Mediterraneo_di_transizione", # row 10
    "Mediterraneo_di_transizione", # row 11
    "Temperato",                   # row 12
    "Temperato",                   # row 13
    "Temperato_di_transizione",    # row 14
    "Temperato_di_transizione",    # row 15
    "Mediterraneo_di_transizione", # row 16
    "Temperato_di_transizione",    # row 17
    "Mediterraneo",                # row 18
    "Mediterraneo_di_transizione", # row 19
    "Temperato_di_transizione",    # row 20
    "Mediterraneo_di_transizione", # row 21
    "Mediterraneo_di_transizione", # row 22
    "Mediterraneo_di_transizione", # row 23
    "Temperato",                   # row 24
    "Mediterraneo_di_transizione", # row 25
    "Mediterraneo",                # row 26
    "Temperato_di_transizione"     # row 27
)

for ($i = 0; $i -lt $climaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $climaValues[$i]
}

# Set column L width to fit the longest value, matching the other bestFit columns
$ws.Columns.Item(12).ColumnWidth = 26.6

# Update the selection to the new column as in the saved file
$ws.Range("L1:L27").Select() | Out-Null
